$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ihme_tab_1.2")

# Row -> (deaths, month, day)
$updates = @(
    @{ Row = 2;  Deaths = 8534.200000000001; Month = "08"; Day = "04" },
    @{ Row = 3;  Deaths = 43962.01;          Month = "08"; Day = "04" },
    @{ Row = 4;  Deaths = 35229.87;          Month = "08"; Day = "04" },
    @{ Row = 5;  Deaths = 27006.12;          Month = "08"; Day = "04" },
    @{ Row = 6;  Deaths = 11093.42;          Month = "08"; Day = "04" },
    @{ Row = 7;  Deaths = 34602.38;          Month = "08"; Day = "04" },
    @{ Row = 8;  Deaths = 9500.376;          Month = "08"; Day = "04" },
    @{ Row = 9;  Deaths = 6598.558;          Month = "08"; Day = "04" },
    @{ Row = 10; Deaths = 3398.807;          Month = "08"; Day = "04" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 5).Value = $u.Deaths

    # Month/day are stored as zero-padded text ("08", "04"), not numbers,
    # so force text formatting before assigning or Excel will coerce them
    # to numeric values and drop the leading zero.
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $u.Month
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $u.Day
}
